$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table: triangle sides (A) and classification (B)
$data = @(
    @("3,4,5", "Scalene"),
    @("0,4,5", "Not a Triangle"),
    @("200,4,5", "Not a Triangle"),
    @("3,0,5", "Not a Triangle"),
    @("3,200,5", "Not a Triangle"),
    @("3,4,0", "Not a Triangle"),
    @("3,4,200", "Not a Triangle"),
    @("1,2,3", "Not a Triangle"),
    @("3,1,2", "Not a Triangle"),
    @("1,3,2", "Not a Triangle"),
    @("3,3,3", "Equilateral"),
    @("2,2,3", "Isosceles"),
    @("3,2,2", "Isosceles"),
    @("2,3,2", "Isosceles")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("C9").Select()
